$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text, matching the source data type
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.896.27'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.648.27'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.81'
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.36'
$ws.Range('E6').Value = '  +1.05%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.59'
$ws.Range('E9').Value = '  +1.44%  '
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.382'
$ws.Range('E11').Value = '  +4.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.156'
$ws.Range('E12').Value = '  +0.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.118.38'
$ws.Range('E13').Value = '  +1.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.04'
$ws.Range('E14').Value = '  +11.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.893.57'
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('E16').Value = '  +0.84%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.660.25'
$ws.Range('E17').Value = '  +1.49%  '
$ws.Range('E19').Value = '  +1.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '349.28'
$ws.Range('E20').Value = '  +0.17%  '
$ws.Range('E21').Value = '  +0.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.997'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('E23').Value = '  +2.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.09'
$ws.Range('E24').Value = '  +1.50%  '
$ws.Range('E25').Value = '  +0.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.996'
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('E27').Value = '  +5.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.98'
$ws.Range('E28').Value = '  +7.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0810'
$ws.Range('E29').Value = '  +2.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.83'
$ws.Range('E30').Value = '  +7.61%  '
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '166.93'
$ws.Range('E32').Value = '  +2.90%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.89'
$ws.Range('E33').Value = '  +1.76%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.53'
$ws.Range('E34').Value = '  +7.15%  '
$ws.Range('E35').Value = '  +7.54%  '
$ws.Range('E36').Value = '  +7.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.66'
$ws.Range('E37').Value = '  +3.76%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '338.37'
$ws.Range('E38').Value = '  +12.67%  '
$ws.Range('E39').Value = '  +4.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.899'
$ws.Range('E40').Value = '  +6.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '38.55'
$ws.Range('E41').Value = '  +1.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.25'
$ws.Range('E42').Value = '  +4.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.34'
$ws.Range('E43').Value = '  +2.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '134.07'
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('E45').Value = '  +1.23%  '
$ws.Range('E46').Value = '  +2.67%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '20.58'
$ws.Range('E47').Value = '  +3.07%  '
$ws.Range('E48').Value = '  +2.58%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.614'
$ws.Range('E49').Value = '  +1.39%  '
$ws.Range('E50').Value = '  +0.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.097.52'
$ws.Range('E51').Value = '  +3.72%  '
